# Update the LR-pair sheet with newly-computed TPM-based values.
# The underlying clustering changed (Neutrophils -> ECs / Resolving-Mac,
# and a MuSCs row replaces the Rspo2->Rspo2 self row order), and the data
# shrinks from 8 target-cluster rows (2 sending clusters x 4 targets) down
# to a single sending cluster (FAPs) x 5 targets.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1681403333333333
$ws.Range("H2").Value = 0.504421
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.184159
$ws.Range("N2").Value = 0.5524770000000001
$ws.Range("O2").Value = 0.1331158254681294
$ws.Range("P2").Value = 0.1331158254681294
$ws.Range("Q2").Value = 0.03096455564633334
$ws.Range("R2").Value = 0.278681000817
$ws.Range("S2").Value = 0.1331158254681294
$ws.Range("T2").Value = 0.1331158254681294

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1681403333333333
$ws.Range("H3").Value = 0.504421
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6233136666666667
$ws.Range("N3").Value = 1.869941
$ws.Range("O3").Value = 0.4505504116763221
$ws.Range("P3").Value = 0.450550411676322
$ws.Range("Q3").Value = 0.1048041676845556
$ws.Range("R3").Value = 0.9432375091609999
$ws.Range("S3").Value = 0.4505504116763221
$ws.Range("T3").Value = 0.450550411676322

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1681403333333333
$ws.Range("H4").Value = 0.504421
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.064319
$ws.Range("N4").Value = 0.192957
$ws.Range("O4").Value = 0.04649176406412185
$ws.Range("P4").Value = 0.04649176406412184
$ws.Range("Q4").Value = 0.01081461809966667
$ws.Range("R4").Value = 0.097331562897
$ws.Range("S4").Value = 0.04649176406412185
$ws.Range("T4").Value = 0.04649176406412184

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1681403333333333
$ws.Range("H5").Value = 0.504421
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4910533333333333
$ws.Range("N5").Value = 1.47316
$ws.Range("O5").Value = 0.3549485488927676
$ws.Range("P5").Value = 0.3549485488927676
$ws.Range("Q5").Value = 0.08256587115111111
$ws.Range("R5").Value = 0.74309284036
$ws.Range("S5").Value = 0.3549485488927676
$ws.Range("T5").Value = 0.3549485488927676

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rspo2"
$ws.Range("C6").Value = "Lgr5"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1681403333333333
$ws.Range("H6").Value = 0.504421
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02060433333333333
$ws.Range("N6").Value = 0.061813
$ws.Range("O6").Value = 0.0148934498986591
$ws.Range("P6").Value = 0.0148934498986591
$ws.Range("Q6").Value = 0.003464419474777778
$ws.Range("R6").Value = 0.031179775273
$ws.Range("S6").Value = 0.0148934498986591
$ws.Range("T6").Value = 0.0148934498986591

# Remove now-unused rows 7-9 so the used range shrinks to A1:T6
$ws.Rows("7:9").Delete()
